$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value as literal text, preserving the cell's default (unstyled) appearance.
# Cryptocurrency prices such as "227.73" would otherwise be auto-detected by Excel
# as numbers; the source data keeps every Price/Volume cell as plain text.
function Set-TextValue($cellAddr, $val) {
    $c = $ws.Range($cellAddr)
    $c.NumberFormat = "@"
    $c.Value2 = $val
    $c.ClearFormats()
}

Set-TextValue "D2" "27.851.14"
$ws.Range("E2").Value2 = "  +6.65%  "
Set-TextValue "D3" "1.741.22"
$ws.Range("E3").Value2 = "  +5.32%  "
$ws.Range("E4").Value2 = "  +0.06%  "
Set-TextValue "D5" "227.73"
$ws.Range("E5").Value2 = "  +4.05%  "
Set-TextValue "D6" "0.5461"
$ws.Range("E6").Value2 = "  +3.87%  "
Set-TextValue "D7" "1.005"
$ws.Range("E7").Value2 = "  +0.05%  "
Set-TextValue "D8" "0.2773"
$ws.Range("E8").Value2 = "  +3.85%  "
Set-TextValue "D9" "0.06738"
$ws.Range("E9").Value2 = "  +5.78%  "
Set-TextValue "D10" "21.85"
$ws.Range("E10").Value2 = "  +6.16%  "
Set-TextValue "D11" "0.07779"
$ws.Range("E11").Value2 = "  +1.04%  "
Set-TextValue "D12" "4.699"
$ws.Range("E12").Value2 = "  +2.19%  "
Set-TextValue "D13" "1.750.20"
$ws.Range("E13").Value2 = "  +5.01%  "
Set-TextValue "D14" "1.983.66"
$ws.Range("E14").Value2 = "  +5.46%  "
Set-TextValue "D15" "0.5976"
$ws.Range("E15").Value2 = "  +6.53%  "
Set-TextValue "D16" "0.0₅8394"
$ws.Range("E16").Value2 = "  +1.84%  "
Set-TextValue "D17" "68.92"
$ws.Range("E17").Value2 = "  +5.34%  "
Set-TextValue "D18" "27.857.86"
$ws.Range("E18").Value2 = "  +6.70%  "
Set-TextValue "D19" "224.23"
$ws.Range("E19").Value2 = "  +17.19%  "
Set-TextValue "D20" "4.843"
$ws.Range("E20").Value2 = "  +3.06%  "
Set-TextValue "D21" "1.003"
$ws.Range("E21").Value2 = "  -0.11%  "
Set-TextValue "D22" "10.93"
$ws.Range("E22").Value2 = "  +5.30%  "
Set-TextValue "D23" "6.242"
$ws.Range("E23").Value2 = "  +4.32%  "
Set-TextValue "D24" "1.005"
$ws.Range("E24").Value2 = "  +0.04%  "
Set-TextValue "D25" "146.30"
$ws.Range("E25").Value2 = "  +0.25%  "
Set-TextValue "D26" "0.1251"
$ws.Range("E26").Value2 = "  +4.08%  "
Set-TextValue "D27" "1.683"
$ws.Range("E27").Value2 = "  +12.47%  "
Set-TextValue "D28" "7.461"
$ws.Range("E28").Value2 = "  +2.73%  "
Set-TextValue "D29" "17.20"
$ws.Range("E29").Value2 = "  +7.78%  "
Set-TextValue "D30" "0.05674"
$ws.Range("E30").Value2 = "  +0.35%  "
Set-TextValue "D31" "1.312"
$ws.Range("E31").Value2 = "  +3.09%  "
Set-TextValue "D32" "3.699"
$ws.Range("E32").Value2 = "  +5.63%  "
Set-TextValue "D33" "3.522"
$ws.Range("E33").Value2 = "  +4.03%  "
Set-TextValue "D34" "1.687"
$ws.Range("E34").Value2 = "  +6.72%  "
Set-TextValue "D35" "0.9775"
$ws.Range("E35").Value2 = "  +3.19%  "
Set-TextValue "D36" "2.856"
$ws.Range("E36").Value2 = "  +1.97%  "
Set-TextValue "D37" "2.450"
$ws.Range("E37").Value2 = "  +1.84%  "
Set-TextValue "D38" "0.5975"
$ws.Range("E38").Value2 = "  +3.21%  "
Set-TextValue "D39" "0.01666"
$ws.Range("E39").Value2 = "  +4.49%  "
Set-TextValue "D40" "5.989"
$ws.Range("E40").Value2 = "  +0.15%  "
Set-TextValue "D41" "0.8516"
$ws.Range("E41").Value2 = "  +1.25%  "
Set-TextValue "D42" "1.046.69"
$ws.Range("E42").Value2 = "  +2.05%  "
Set-TextValue "D43" "1.004"
$ws.Range("E43").Value2 = "  +0.08%  "
Set-TextValue "D44" "102.03"
$ws.Range("E44").Value2 = "  +0.34%  "
Set-TextValue "D45" "1.888.71"
$ws.Range("E45").Value2 = "  +5.38%  "
$ws.Range("B46").Value2 = "Aave"
$ws.Range("C46").Value2 = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
Set-TextValue "D46" "59.62"
$ws.Range("E46").Value2 = "  +1.77%  "
$ws.Range("B47").Value2 = "BabyDogeCoin"
$ws.Range("C47").Value2 = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
Set-TextValue "D47" "0.0₈110"
$ws.Range("E47").Value2 = "  +7.32%  "
Set-TextValue "D48" "8.303"
$ws.Range("E48").Value2 = "  +3.08%  "
Set-TextValue "D49" "0.4441"
$ws.Range("E49").Value2 = "  +2.26%  "
$ws.Range("B50").Value2 = "Frax"
$ws.Range("C50").Value2 = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
Set-TextValue "D50" "1.005"
$ws.Range("E50").Value2 = "  +0.00%  "
$ws.Range("B51").Value2 = "Cronos"
$ws.Range("C51").Value2 = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
Set-TextValue "D51" "0.05327"
$ws.Range("E51").Value2 = "  -0.16%  "
